$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.310.22'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.929.99'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7487'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3176'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.52'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07117'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7808'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08052'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').Value = '1.971.82'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.405'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '30.318.77'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.049'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007933'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.136.32'
$ws.Range('E22').Value = '  -2.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.661'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.589'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.34'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1293'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.188'
$ws.Range('D29').Style = 'Normal'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.566'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.38%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.362'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.428'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.142'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05252'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.320'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7578'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.766'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01954'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.799'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.522'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '77.40'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4529'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.971'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8442'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.714'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').Value = '2.109.36'
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1223'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.18%  '
